$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '97.293.74'
$ws.Range("E2").Value = '  -0.67%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.354.27'
$ws.Range("E3").Value = '  -1.73%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.58'
$ws.Range("E5").Value = '  -1.68%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '659.60'
$ws.Range("E6").Value = '  +0.46%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.45'
$ws.Range("E7").Value = '  -2.67%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.428'
$ws.Range("E8").Value = '  -2.93%  '

$ws.Range("E9").Value = '  -0.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.03'
$ws.Range("E10").Value = '  -5.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.352.03'
$ws.Range("E11").Value = '  -1.63%  '

$ws.Range("E12").Value = '  -2.61%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.52'
$ws.Range("E13").Value = '  -1.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '96.930.92'
$ws.Range("E14").Value = '  -0.76%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.13'
$ws.Range("E15").Value = '  -5.40%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000256'
$ws.Range("E16").Value = '  -2.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.985.30'
$ws.Range("E17").Value = '  -1.96%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.94'
$ws.Range("E18").Value = '  +3.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.318.18'
$ws.Range("E19").Value = '  -3.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.582'
$ws.Range("E20").Value = '  +17.59%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.89'
$ws.Range("E21").Value = '  +1.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.78'
$ws.Range("E22").Value = '  +0.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '513.24'
$ws.Range("E23").Value = '  +0.33%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.38'
$ws.Range("E24").Value = '  -3.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000201'
$ws.Range("E25").Value = '  -3.26%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.68'
$ws.Range("E26").Value = '  +7.18%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '97.41'
$ws.Range("E27").Value = '  -1.82%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.29'
$ws.Range("E28").Value = '  -5.47%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.149'
$ws.Range("E29").Value = '  -2.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '11.54'
$ws.Range("E30").Value = '  +0.12%  '

$ws.Range("E31").Value = '  +0.00%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.190'
$ws.Range("E32").Value = '  -5.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.55'
$ws.Range("E33").Value = '  +11.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.563'
$ws.Range("E35").Value = '  -1.95%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '28.63'
$ws.Range("E36").Value = '  -3.97%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.52'
$ws.Range("E37").Value = '  +5.27%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.95'
$ws.Range("E38").Value = '  +0.91%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.153'
$ws.Range("E39").Value = '  -0.68%  '

$ws.Range("E40").Value = '  +0.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '512.81'
$ws.Range("E41").Value = '  -3.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0437'
$ws.Range("E42").Value = '  +3.84%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '24.40'
$ws.Range("E43").Value = '  -1.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.849'
$ws.Range("E44").Value = '  -2.14%  '

$ws.Range("B45").Value = 'ImmutableX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.70'
$ws.Range("E45").Value = '  +6.05%  '

$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.66'
$ws.Range("E46").Value = '  +2.32%  '

$ws.Range("B47").Value = 'MantraDAO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.66'
$ws.Range("E47").Value = '  -1.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.68'
$ws.Range("E48").Value = '  +4.22%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.56'
$ws.Range("E49").Value = '  +6.35%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.15'
$ws.Range("E50").Value = '  -5.50%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.02'
$ws.Range("E51").Value = '  -3.28%  '
